$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mississippi (row 42) failed to fetch data on this run -- clear the measurement
# columns (B..H) to an empty text value, matching the "request failed" pattern
# already used elsewhere in this sheet (e.g. row 19, California).
foreach ($col in 2..8) {
    $cell = $ws.Cells.Item(42, $col)
    $cell.Value = "'"
    $cell.Style = "Normal"
}

# Booleans flip to False since no data was retrieved this run.
$ws.Cells.Item(42, 9).Value = $false
$ws.Cells.Item(42, 10).Value = $false

# Status message reflects the connection error encountered during this run.
$ws.Cells.Item(42, 15).Value = "An error occurred. ... ConnectionError(ProtocolError('Connection aborted.', ConnectionResetError(104, 'Connection reset by peer')))"
